{"js": "// The TODO list currently reads:\n//   1. Layout for subs index and file search to match manage users index\n//   2. Email templates\n//   3. Registration to add option to register as an mga or uw\n//   4. Queue sending emails\n//   5. Send thank you email and generate random sub # in that email - finish\n//   6. Attachments in emails\n//   7. Factories for rw\n//   8. Factories for quotes\n//\n// It needs to become:\n//   1. Email templates\n//   2. Queue sending emails\n//   3. Attachments in emails\n//\n// Locate every affected paragraph up front (by its current, still-unique\n// text) before mutating anything, so none of the edits below interfere\n// with one another.\nconst body = context.document.body;\n\nfunction findParagraph(text) {\n  const results = body.search(text, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  return results;\n}\n\nconst layoutResults = findParagraph(\n  \"Layout for subs index and file search to match manage users index\"\n);\nconst emailTemplatesResults = findParagraph(\"Email templates\");\nconst registrationResults = findParagraph(\n  \"Registration to add option to register as an mga or uw\"\n);\nconst queueResults = findParagraph(\"Queue sending emails\");\nconst thankYouResults = findParagraph(\n  \"Send thank you email and generate random sub # in that email\"\n);\nconst factoriesRwResults = findParagraph(\"Factories for rw\");\nconst factoriesQuotesResults = findParagraph(\"Factories for quotes\");\n\nawait context.sync();\n\nconst layoutParagraph = layoutResults.items[0].paragraphs.getFirst();\nconst emailTemplatesParagraph = emailTemplatesResults.items[0].paragraphs.getFirst();\nconst registrationParagraph = registrationResults.items[0].paragraphs.getFirst();\nconst queueParagraph = queueResults.items[0].paragraphs.getFirst();\nconst thankYouParagraph = thankYouResults.items[0].paragraphs.getFirst();\nconst factoriesRwParagraph = factoriesRwResults.items[0].paragraphs.getFirst();\nconst factoriesQuotesParagraph = factoriesQuotesResults.items[0].paragraphs.getFirst();\n\n// Re-purpose the first two surviving paragraphs with their new text...\nlayoutParagraph.insertText(\"Email templates\", Word.InsertLocation.replace);\nemailTemplatesParagraph.insertText(\"Queue sending emails\", Word.InsertLocation.replace);\n\n// ...and drop the paragraphs that are no longer needed.\nregistrationParagraph.delete();\nqueueParagraph.delete();\nthankYouParagraph.delete();\nfactoriesRwParagraph.delete();\nfactoriesQuotesParagraph.delete();\n\n// \"Attachments in emails\" is left untouched, in place.\n\nawait context.sync();\n", "ps1": "# The TODO list currently reads:\n#   1. Layout for subs index and file search to match manage users index\n#   2. Email templates\n#   3. Registration to add option to register as an mga or uw\n#   4. Queue sending emails\n#   5. Send thank you email and generate random sub # in that email - finish\n#   6. Attachments in emails\n#   7. Factories for rw\n#   8. Factories for quotes\n#\n# It needs to become:\n#   1. Email templates\n#   2. Queue sending emails\n#   3. Attachments in emails\n#\n# Snapshot every paragraph up front (by its current, still-unique text) so\n# that later edits/deletions can't confuse the lookups below.\n$d = $word.ActiveDocument\n$paragraphs = @($d.Paragraphs)\n\n$layoutParagraph = $null\n$emailTemplatesParagraph = $null\n$registrationParagraph = $null\n$queueParagraph = $null\n$thankYouParagraph = $null\n$factoriesRwParagraph = $null\n$factoriesQuotesParagraph = $null\n\nforeach ($p in $paragraphs) {\n    $text = $p.Range.Text\n\n    if ($text -like \"Layout for subs index and file search to match manage users index*\") {\n        $layoutParagraph = $p\n    }\n    elseif ($text -like \"Email templates*\") {\n        $emailTemplatesParagraph = $p\n    }\n    elseif ($text -like \"Registration to add option to register as an mga or uw*\") {\n        $registrationParagraph = $p\n    }\n    elseif ($text -like \"Queue sending emails*\") {\n        $queueParagraph = $p\n    }\n    elseif ($text -like \"Send thank you email and generate random sub # in that email*\") {\n        $thankYouParagraph = $p\n    }\n    elseif ($text -like \"Factories for rw*\") {\n        $factoriesRwParagraph = $p\n    }\n    elseif ($text -like \"Factories for quotes*\") {\n        $factoriesQuotesParagraph = $p\n    }\n}\n\n# Re-purpose the first two surviving paragraphs with their new text...\n$layoutParagraph.Range.Text = \"Email templates\"\n$emailTemplatesParagraph.Range.Text = \"Queue sending emails\"\n\n# ...and drop the paragraphs that are no longer needed. Delete from the\n# bottom of the document upward so that removing one paragraph never\n# invalidates an already-captured reference to a later paragraph.\n$factoriesQuotesParagraph.Range.Delete()\n$factoriesRwParagraph.Range.Delete()\n$thankYouParagraph.Range.Delete()\n$queueParagraph.Range.Delete()\n$registrationParagraph.Range.Delete()\n\n# \"Attachments in emails\" is left untouched, in place.\n"}
